$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update asset names (column A)
$ws.Range("A2").Value = "Elkjele e"
$ws.Range("A3").Value = "Elkjele f"
$ws.Range("A4").Value = "Elkjele g"

# Update Brønnøysund register numbers (column E)
$ws.Range("E2").Value = "'707057500051530000"
$ws.Range("E3").Value = "'707057500053255000"
$ws.Range("E4").Value = "707057500056855200"

# Update MPID values (column H)
$ws.Range("H2").Value = 926956191
$ws.Range("H3").Value = 918874321
$ws.Range("H4").Value = 913303334

# Update the view: scroll / zoom / selection
$ws.Application.ActiveWindow.Zoom = 258
$ws.Range("E7").Select()
